# Update gh-pages output (杭州-漫展信息.xlsx) to the data generated at 456a3b4.
#
# Changes:
#  1. Sheet "展览": bump "想去人数" (F column) counters for a number of
#     existing exhibitions, and insert a brand-new exhibition row
#     (杭州·星部落动漫嘉年华) as the new row 37 -- shifting the three
#     exhibitions that used to be rows 37-39 down to rows 38-40.
#  2. Sheet "本地生活": bump F-column counters for its 3 data rows.
#  3. Sheet "全部类型" (the combined/all-types view): bump the matching
#     F-column counters for the same underlying events.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Simple counter bumps (no structural change for these rows).
$ws1.Range("F3").Value = 8115
$ws1.Range("F4").Value = 1911
$ws1.Range("F5").Value = 6499
$ws1.Range("F7").Value = 2050
$ws1.Range("F8").Value = 565
$ws1.Range("F9").Value = 42
$ws1.Range("F15").Value = 8471
$ws1.Range("F20").Value = 1803
$ws1.Range("F25").Value = 20
$ws1.Range("F30").Value = 2054
$ws1.Range("F31").Value = 842
$ws1.Range("F32").Value = 465
$ws1.Range("F35").Value = 172

# Insert a new row at position 37 (existing rows 37-39 shift down to 38-40).
$ws1.Rows.Item(37).Insert()

# The freshly inserted row has no formatting of its own yet; clone the
# look (number formats / borders / alignment) of the row directly below
# it -- which, after the shift, holds what used to be row 37 -- so the
# new row matches the sheet's established per-column styling.
$ws1.Range("A38:I38").Copy()
$ws1.Range("A37:I37").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Populate the new row with the "星部落动漫嘉年华" exhibition.
$ws1.Range("A37").Value = 36
# Column B stores dates as plain text (e.g. "2024-11-02"), not real Excel
# dates -- force a text number format first so the literal string is kept
# instead of being auto-converted into a date serial number.
$ws1.Range("B37").NumberFormat = "@"
$ws1.Range("B37").Value = "2024-11-02"
$ws1.Range("C37").Value = "杭州·星部落动漫嘉年华"
$ws1.Range("D37").Value = "康候圣街99号 顺丰创新中心"
$ws1.Range("E37").Value = "2024.11.02 09:00-11.03 16:00"
$ws1.Range("F37").Value = 0
$ws1.Range("G37").Value = 49
$ws1.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=91795"
$ws1.Range("I37").Value = "//i1.hdslb.com/bfs/openplatform/202408/KCwYmgHz1724908471827.jpeg"

# The explicit text NumberFormat above leaves column B with a distinct
# style from its siblings; re-sync it from the neighboring (already
# correctly-styled) text cell below so it matches the rest of the column.
$ws1.Range("B38").Copy()
$ws1.Range("B37").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Row 38 now holds what used to be row 37 (杭州·New World动漫博览会):
# its running index and "want to go" count need updating; everything
# else about the event is unchanged.
$ws1.Range("A38").Value = 37
$ws1.Range("F38").Value = 16

# Rows 39 and 40 now hold what used to be rows 38 and 39 respectively --
# their content is unchanged, only the running index needs to follow the
# new row numbers.
$ws1.Range("A39").Value = 38
$ws1.Range("A40").Value = 39

# ---------------------------------------------------------------------
# Sheet "本地生活" (local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2319
$ws3.Range("F3").Value = 707
$ws3.Range("F4").Value = 304

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types combined)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2319
$ws4.Range("F3").Value = 707
$ws4.Range("F6").Value = 8115
$ws4.Range("F8").Value = 304
$ws4.Range("F9").Value = 1911
$ws4.Range("F10").Value = 6499
$ws4.Range("F11").Value = 2050
$ws4.Range("F13").Value = 565
$ws4.Range("F14").Value = 42
$ws4.Range("F23").Value = 8471
$ws4.Range("F28").Value = 1803
$ws4.Range("F35").Value = 2054
$ws4.Range("F36").Value = 842
$ws4.Range("F38").Value = 465
